$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.689.52'
$ws.Range("E2").Value = '  +0.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.691.54'
$ws.Range("E3").Value = '  +0.78%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.59%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.97'
$ws.Range("E5").Value = '  +1.56%  '
$ws.Range("E6").Value = '  +0.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3717'
$ws.Range("E7").Value = '  +0.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.85'
$ws.Range("E8").Value = '  +1.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3397'
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.201'
$ws.Range("E10").Value = '  +3.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07418'
$ws.Range("E11").Value = '  +2.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.003'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.286'
$ws.Range("E13").Value = '  +2.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.82'
$ws.Range("E14").Value = '  +3.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.942'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.692.94'
$ws.Range("E16").Value = '  +1.02%  '
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06692'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '82.78'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.08'
$ws.Range("E21").Value = '  +3.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.286'
$ws.Range("E22").Value = '  +3.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.86'
$ws.Range("E23").Value = '  +6.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.684.71'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.445'
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.748'
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.08'
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '148.28'
$ws.Range("E28").Value = '  -3.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '130.95'
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.882.26'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.217'
$ws.Range("E31").Value = '  +24.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.667'
$ws.Range("E32").Value = '  +6.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.207'
$ws.Range("E33").Value = '  +4.47%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.50'
$ws.Range("E34").Value = '  +9.45%  '
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08659'
$ws.Range("E35").Value = '  +2.69%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.749'
$ws.Range("E36").Value = '  +3.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.496'
$ws.Range("E37").Value = '  +3.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06571'
$ws.Range("E38").Value = '  +3.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.006'
$ws.Range("E39").Value = '  +4.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02388'
$ws.Range("E40").Value = '  +3.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2203'
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.256'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6351'
$ws.Range("E43").Value = '  +4.43%  '
$ws.Range("E44").Value = '  +0.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.66'
$ws.Range("E45").Value = '  +5.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6045'
$ws.Range("E46").Value = '  +2.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.811'
$ws.Range("E47").Value = '  +1.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.097'
$ws.Range("E48").Value = '  +4.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '128.40'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07219'
$ws.Range("E50").Value = '  +0.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.88'
$ws.Range("E51").Value = '  +4.26%  '
